# "Sep 2020 to Oct 2020" monthly Cord Cutting Wizard change-log upload:
# the sheet that tracked "Jul 2020 to Aug 2020" is renamed for the new
# reporting period. Excel automatically re-points the sheet-scoped
# _FilterDatabase defined name (and any formulas) at the new sheet name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sep 2020 to Oct 2020"
